# Commit: Specify log2 when referring to fold change (Closes #2)

$wb = $excel.ActiveWorkbook

# --- 1. "interactions" sheet: rename fold-change headers to l2fc ---
$wsInter = $wb.Worksheets.Item("interactions")
$wsInter.Range("D1").Value = "opc_l2fc_a"
$wsInter.Range("E1").Value = "endo_l2fc_a"
$wsInter.Range("F1").Value = "dac_l2fc_a"
$wsInter.Range("G1").Value = "opc_l2fc_b"
$wsInter.Range("H1").Value = "endo_l2fc_b"
$wsInter.Range("I1").Value = "dac_l2fc_b"

# --- 2. "gene_counts" sheet: re-order gene rows (counts recomputed / re-sorted) ---
$wsCounts = $wb.Worksheets.Item("gene_counts")

$genes = @("Itgb1","Ptprc","Il6","Cxcl12","Itga4","Vcam1","Anxa1","B2m","Itga7","Il15","Cxcr4","Il7r","Itga5","Itgam","Kit","Alcam","Itgav","Cd1d1","Nt5e","Itga6","Tek","Selp","Dpp4")
$counts = @(20,12,10,8,8,6,6,6,4,4,4,4,4,4,4,2,2,2,2,2,2,2,2)

for ($i = 0; $i -lt $genes.Length; $i++) {
    $row = $i + 2
    $wsCounts.Cells.Item($row, 1).Value = $genes[$i]
    $wsCounts.Cells.Item($row, 2).Value = $counts[$i]
}

# --- 3. "parameters" sheet: rename cutoff, insert new ratio parameter, shift others down ---
$wsParams = $wb.Worksheets.Item("parameters")

$wsParams.Range("A2").Value = "log2_fold_change_cutoff"
$wsParams.Range("B2").Value = 1.5

$wsParams.Range("A3").Value = "opc_microglia_min_ratio"
$wsParams.Range("B3").Value = 0.05

$wsParams.Range("A4").Value = "string_score_cutoff"
$wsParams.Range("B4").Value = 700

$wsParams.Range("A5").Value = "GO_terms"
$wsParams.Range("B5").Value = "['cell surface', 'cell-cell adhesion']"
